# Applies the Thu Feb 16 13:38:06 UTC 2023 symbol-list refresh:
# updates the "Price" (column D) and "Volume(1h)" (column E) text values
# for the rows whose quoted figures moved in this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new display text (values are stored as plain
# text in this sheet, e.g. "319.34" / "5.45%", so each one is entered with a
# leading apostrophe to keep Excel from re-typing it as a number/percentage,
# then the style is put back to "Normal" so no formatting changes stick.
$updates = [ordered]@{
    "D2" = "319.34"
    "E2" = "5.45%"
    "D3" = "48.62"
    "E3" = "12.61%"
    "D4" = "5.271"
    "E4" = "4.74%"
    "D5" = "0.08048"
    "E5" = "4.87%"
    "D6" = "4.585"
    "E6" = "4.13%"
    "D7" = "1.421"
    "E7" = "35.25%"
    "D8" = "1.648"
    "E8" = "2.58%"
    "D9" = "0.1283"
    "E9" = "4.42%"
    "D10" = "0.1942"
    "E10" = "4.88%"
    "D11" = "0.09374"
    "E11" = "3.80%"
    "D12" = "0.04580"
    "E12" = "9.86%"
    "D13" = "0.1046"
    "E13" = "-0.01%"
    "D14" = "0.001337"
    "E14" = "5.52%"
    "E15" = "0.33%"
    "D16" = "0.005851"
    "E16" = "1.31%"
    "D17" = "3.332"
    "E17" = "0.46%"
    "D18" = "2.439"
    "E18" = "2.32%"
    "D19" = "0.3419"
    "E19" = "2.36%"
    "D20" = "8.178"
    "E20" = "-2.80%"
    "D21" = "0.1401"
    "E21" = "-0.17%"
    "D22" = "0.3099"
    "E22" = "6.99%"
    "D23" = "0.001304"
    "E23" = "2.37%"
    "D24" = "0.004245"
    "E24" = "-5.36%"
    "D25" = "0.0001348"
    "E25" = "-0.05%"
    "D26" = "0.0003541"
    "E26" = "-95.23%"
    "D38" = "0.02689"
    "E38" = "10.00%"
    "D39" = "0.05706"
    "E39" = "8.21%"
    "E40" = "4.13%"
    "D41" = "0.007978"
    "E41" = "3.49%"
    "D42" = "0.1443"
    "E42" = "7.02%"
    "D43" = "0.007685"
    "E43" = "4.53%"
    "D44" = "0.007868"
    "E44" = "-7.12%"
    "D45" = "0.3501"
    "E45" = "15.02%"
    "D46" = "0.00006899"
    "E46" = "4.06%"
    "E47" = "0.14%"
    "D48" = "0.05494"
    "E48" = "19.80%"
    "D49" = "0.004002"
    "E49" = "-4.74%"
    "D50" = "0.00002101"
    "E50" = "0.14%"
    "D51" = "0.0002001"
    "E51" = "0.14%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $updates[$ref]
    $cell.Style = "Normal"
}
